# Update crypto price/volume figures per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.262.84"
$ws.Range("E2").Value = "  +2.94%  "
$ws.Range("D3").Value = "2.315.72"
$ws.Range("E3").Value = "  +2.60%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.66"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.50"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.78%  "
$ws.Range("E7").Value = "  +2.13%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +7.28%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.08"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.66%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0814"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.24%  "
$ws.Range("E12").Value = "  +0.76%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.03"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.84%  "
$ws.Range("D14").Value = "2.673.33"
$ws.Range("E14").Value = "  +2.54%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.99"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.32%  "
$ws.Range("D16").Value = "2.316.72"
$ws.Range("E16").Value = "  +2.15%  "
$ws.Range("E17").Value = "  +3.01%  "
$ws.Range("D18").Value = "43.171.83"
$ws.Range("E18").Value = "  +2.98%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.55"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.76%  "
$ws.Range("D20").Value = "0.0₃0920"
$ws.Range("E20").Value = "  +2.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.14"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.38%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.46"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "241.50"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.95%  "
$ws.Range("E24").Value = "  +6.47%  "
$ws.Range("E25").Value = "  +3.21%  "
$ws.Range("E27").Value = "  +5.18%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "37.46"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.99%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.66"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.49%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.11"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.10%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "167.95"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.46%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.33"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.66%  "
$ws.Range("E33").Value = "  +0.05%  "
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "17.96"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.61%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0744"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.40%  "
$ws.Range("E37").Value = "  +3.06%  "
$ws.Range("E38").Value = "  +1.27%  "
$ws.Range("E39").Value = "  +3.02%  "
$ws.Range("E40").Value = "  +2.38%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.28"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +7.76%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "19.82"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.20%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.30"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.23%  "
$ws.Range("E44").Value = "  +3.88%  "
$ws.Range("D45").Value = "1.974.80"
$ws.Range("E45").Value = "  +0.83%  "
$ws.Range("E46").Value = "  +4.41%  "
$ws.Range("E47").Value = "  -0.82%  "
$ws.Range("E48").Value = "  +18.40%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "55.84"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.21%  "
$ws.Range("D50").Value = "2.540.56"
$ws.Range("E50").Value = "  +2.44%  "
$ws.Range("E51").Value = "  +4.46%  "
